# Generate Report for Handback
# This records a handback for the "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md"
# file (row 7 of both the "zh-cn" and "de-de" sheets), and flags that the
# handback file version is not the latest with an explanatory Error Detail.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d509be54024f65374340eedc3185aaeed01b8f07/e2e/cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcad4860e6c31287cd8c0ccd0427ca494391260f/e2e/cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md."

$hyperlinkColor = 15570276

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# I7: Latest Target File
$wsZh.Range("I7").Value = "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bcad4860e6c31287cd8c0ccd0427ca494391260f/e2e/cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md", "", "", "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md")
$wsZh.Range("I7").Font.Underline = $true
$wsZh.Range("I7").Font.Color = $hyperlinkColor

# J7: Latest Handback File
$wsZh.Range("J7").Value = "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.8d17e9b0d629da752e0912f8f057624a4d65c319.zh-cn.xlf"

# K7: Latest Handback DateTime
$wsZh.Range("K7").Value = "2016-08-27 00:55:33"

# P7: Error Detail
$wsZh.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

# I7: Latest Target File
$wsDe.Range("I7").Value = "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bcad4860e6c31287cd8c0ccd0427ca494391260f/e2e/cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md", "", "", "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.md")
$wsDe.Range("I7").Font.Underline = $true
$wsDe.Range("I7").Font.Color = $hyperlinkColor

# J7: Latest Handback File
$wsDe.Range("J7").Value = "cb6b3aa8-f4a0-4a3e-a90c-bdc33a054bae.8d17e9b0d629da752e0912f8f057624a4d65c319.de-de.xlf"

# K7: Latest Handback DateTime
$wsDe.Range("K7").Value = "2016-08-27 00:55:39"

# P7: Error Detail (same message text as zh-cn)
$wsDe.Range("P7").Value = $errorDetail
